# Update the yearly dollar income-statement sheet:
#  - roll the 5-year window forward one year (drop 1396/12, add 1401/12)
#  - roll the "publish date" row forward to match, with the new
#    read_price-algorithm date format for the newest column
#  - refresh every financial figure with the latest database values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: financial period headers (D:H) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish-date headers (D:H) ---
$ws.Range("D9").Value = "1399-05-09 (11)"
$ws.Range("E9").Value = "1400-05-07 (12)"
$ws.Range("F9").Value = "1401-05-09 (9)"
$ws.Range("G9").Value = "1402-02-30 (8)"
$ws.Range("H9").Value = "1402-02-30 (2)"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 2320718
$ws.Range("E11").Value = 3051549
$ws.Range("F11").Value = 3423389
$ws.Range("G11").Value = 5519228
$ws.Range("H11").Value = 4553930

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -1233755
$ws.Range("E12").Value = -1757524
$ws.Range("F12").Value = -1796983
$ws.Range("G12").Value = -2660377
$ws.Range("H12").Value = -2937157

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 1086963
$ws.Range("E13").Value = 1294025
$ws.Range("F13").Value = 1626406
$ws.Range("G13").Value = 2858851
$ws.Range("H13").Value = 1616773

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -76045
$ws.Range("E14").Value = -93062
$ws.Range("F14").Value = -82614
$ws.Range("G14").Value = -130102
$ws.Range("H14").Value = -151059

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense) ---
$ws.Range("D16").Value = 168816
$ws.Range("E16").Value = -11501
$ws.Range("F16").Value = 28112
$ws.Range("G16").Value = 46444
$ws.Range("H16").Value = 66091

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 1179734
$ws.Range("E17").Value = 1189462
$ws.Range("F17").Value = 1571904
$ws.Range("G17").Value = 2775193
$ws.Range("H17").Value = 1531805

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = -101365
$ws.Range("E18").Value = -114673
$ws.Range("F18").Value = -133759
$ws.Range("G18").Value = -94544
$ws.Range("H18").Value = -105363

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense) ---
$ws.Range("D19").Value = 228203
$ws.Range("E19").Value = 242608
$ws.Range("F19").Value = 398093
$ws.Range("G19").Value = 699929
$ws.Range("H19").Value = 285879

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit from continuing ops) ---
$ws.Range("D20").Value = 1306573
$ws.Range("E20").Value = 1317398
$ws.Range("F20").Value = 1836238
$ws.Range("G20").Value = 3380579
$ws.Range("H20").Value = 1712321

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -97512
$ws.Range("E21").Value = -164003
$ws.Range("F21").Value = -143311
$ws.Range("G21").Value = -193981
$ws.Range("H21").Value = -161403

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops) ---
$ws.Range("D22").Value = 1209061
$ws.Range("E22").Value = 1153395
$ws.Range("F22").Value = 1692927
$ws.Range("G22").Value = 3186597
$ws.Range("H22").Value = 1550918

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 1209061
$ws.Range("E24").Value = 1153395
$ws.Range("F24").Value = 1692927
$ws.Range("G24").Value = 3186597
$ws.Range("H24").Value = 1550918

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 2066092
$ws.Range("E26").Value = 2284030
$ws.Range("F26").Value = 1295873
$ws.Range("G26").Value = 1110466
$ws.Range("H26").Value = 2266992
